# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '31.228.95'
$ws.Range('E2').Value = '  +2.34%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.979.88'
$ws.Range('E3').Value = '  +4.92%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9981'
$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7871'
$ws.Range('E5').Value = '  +66.73%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '252.89'
$ws.Range('E6').Value = '  +3.57%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9985'
$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3379'
$ws.Range('E8').Value = '  +16.77%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.50'
$ws.Range('E9').Value = '  +14.51%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06900'
$ws.Range('E10').Value = '  +6.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8490'
$ws.Range('E11').Value = '  +16.80%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08132'
$ws.Range('E12').Value = '  +4.80%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '101.99'
$ws.Range('E13').Value = '  +6.44%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.980.64'
$ws.Range('E14').Value = '  +4.92%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.521'
$ws.Range('E15').Value = '  +6.45%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '276.12'
$ws.Range('E16').Value = '  -2.08%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.229.58'
$ws.Range('E17').Value = '  +2.36%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.98'
$ws.Range('E18').Value = '  +6.90%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007842'
$ws.Range('E19').Value = '  +4.95%  '

$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.234.86'
$ws.Range('E20').Value = '  +4.53%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.672'
$ws.Range('E21').Value = '  +7.61%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.04%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9977'
$ws.Range('E23').Value = '  -0.24%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.804'
$ws.Range('E24').Value = '  +7.77%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.619'
$ws.Range('E25').Value = '  +5.98%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1538'
$ws.Range('E26').Value = '  +59.24%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.28'
$ws.Range('E27').Value = '  +0.68%  '

$ws.Range('E28').Value = '  +3.80%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.215'
$ws.Range('E29').Value = '  +17.04%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.563'
$ws.Range('E30').Value = '  +6.16%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.347'
$ws.Range('E31').Value = '  +0.94%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.544'
$ws.Range('E32').Value = '  +6.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.337'
$ws.Range('E33').Value = '  +4.58%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05150'
$ws.Range('E34').Value = '  +6.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.222'
$ws.Range('E35').Value = '  +8.62%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7408'
$ws.Range('E36').Value = '  +6.86%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.789'
$ws.Range('E37').Value = '  +2.71%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9973'
$ws.Range('E38').Value = '  -0.19%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01982'
$ws.Range('E39').Value = '  +5.31%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.900'
$ws.Range('E40').Value = '  +2.70%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.611'
$ws.Range('E41').Value = '  +6.44%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '78.62'
$ws.Range('E42').Value = '  +5.05%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4647'
$ws.Range('E43').Value = '  +8.82%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.077'
$ws.Range('E44').Value = '  +5.71%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '105.63'
$ws.Range('E45').Value = '  +4.52%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.8506'
$ws.Range('E46').Value = '  +2.97%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9986'
$ws.Range('E47').Value = '  -0.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.04'
$ws.Range('E48').Value = '  +4.13%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.499'
$ws.Range('E49').Value = '  +7.72%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.41'
$ws.Range('E50').Value = '  +3.50%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4259'
$ws.Range('E51').Value = '  +8.24%  '
